$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-12-04 Thursday" "2025-12-05 Friday"

Replace-Text "31×41=1271" "21×84=1764"
Replace-Text "84×38=3192" "47×30=1410"
Replace-Text "68×49=3332" "45×49=2205"
Replace-Text "45×83=3735" "47×51=2397"
Replace-Text "59×35=2065" "97×56=5432"

Replace-Text "18×82=1476" "55×36=1980"
Replace-Text "20×59=1180" "40×65=2600"
Replace-Text "29×70=2030" "16×33=528"
Replace-Text "53×38=2014" "54×83=4482"
Replace-Text "22×19=418" "28×79=2212"

Replace-Text "87×82=7134" "14×60=840"
Replace-Text "29×93=2697" "51×69=3519"
Replace-Text "77×37=2849" "37×25=925"
Replace-Text "67×13=871" "87×24=2088"
Replace-Text "66×61=4026" "28×38=1064"

Replace-Text "73×21=1533" "28×45=1260"
Replace-Text "70×28=1960" "29×32=928"
Replace-Text "50×12=600" "58×26=1508"
Replace-Text "19×85=1615" "33×34=1122"
Replace-Text "97×57=5529" "27×17=459"

Replace-Text "77×61=4697" "88×95=8360"
Replace-Text "40×81=3240" "86×51=4386"
Replace-Text "13×38=494" "64×71=4544"
Replace-Text "73×57=4161" "77×20=1540"
Replace-Text "86×92=7912" "83×94=7802"
